$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '68.387.10'
Set-TextValue 'E2' '  -1.55%  '

Set-TextValue 'D3' '2.448.75'
Set-TextValue 'E3' '  -1.65%  '

Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.11%  '

Set-TextValue 'D5' '554.55'
Set-TextValue 'E5' '  -2.37%  '

Set-TextValue 'D6' '161.00'
Set-TextValue 'E6' '  -1.86%  '

Set-TextValue 'E7' '  -0.11%  '

Set-TextValue 'D8' '0.499'
Set-TextValue 'E8' '  -2.49%  '

Set-TextValue 'D9' '2.448.86'
Set-TextValue 'E9' '  -1.58%  '

Set-TextValue 'D10' '0.148'
Set-TextValue 'E10' '  -6.54%  '

Set-TextValue 'E11' '  -1.48%  '

Set-TextValue 'D12' '0.334'
Set-TextValue 'E12' '  -5.53%  '

Set-TextValue 'D13' '4.73'
Set-TextValue 'E13' '  -3.14%  '

Set-TextValue 'D14' '2.892.67'
Set-TextValue 'E14' '  -1.79%  '

Set-TextValue 'D15' '68.115.97'
Set-TextValue 'E15' '  -1.72%  '

Set-TextValue 'D16' '0.0000166'
Set-TextValue 'E16' '  -4.82%  '

Set-TextValue 'D17' '23.15'
Set-TextValue 'E17' '  -4.64%  '

Set-TextValue 'D18' '2.441.56'
Set-TextValue 'E18' '  -2.00%  '

Set-TextValue 'D19' '10.73'
Set-TextValue 'E19' '  -4.05%  '

Set-TextValue 'D20' '339.89'
Set-TextValue 'E20' '  -1.33%  '

Set-TextValue 'D21' '6.99'
Set-TextValue 'E21' '  -4.93%  '

Set-TextValue 'D22' '3.74'
Set-TextValue 'E22' '  -3.05%  '

Set-TextValue 'E23' '  -0.03%  '

Set-TextValue 'D24' '1.87'
Set-TextValue 'E24' '  -2.70%  '

Set-TextValue 'D25' '66.15'
Set-TextValue 'E25' '  -4.93%  '

Set-TextValue 'D26' '3.66'
Set-TextValue 'E26' '  -5.82%  '

Set-TextValue 'D27' '2.570.31'
Set-TextValue 'E27' '  -1.81%  '

Set-TextValue 'D28' '0.994'
Set-TextValue 'E28' '  -0.16%  '

Set-TextValue 'D29' '8.02'
Set-TextValue 'E29' '  -7.25%  '

Set-TextValue 'D30' '0.0₃0813'
Set-TextValue 'E30' '  -6.72%  '

Set-TextValue 'D31' '7.09'
Set-TextValue 'E31' '  -7.59%  '

Set-TextValue 'E32' '  -0.17%  '

Set-TextValue 'D33' '425.59'
Set-TextValue 'E33' '  -3.67%  '

Set-TextValue 'D34' '1.13'
Set-TextValue 'E34' '  -4.30%  '

Set-TextValue 'D35' '1.61'
Set-TextValue 'E35' '  -5.57%  '

Set-TextValue 'D36' '157.26'
Set-TextValue 'E36' '  +1.64%  '

Set-TextValue 'D37' '18.99'
Set-TextValue 'E37' '  -0.32%  '

Set-TextValue 'E38' '  +0.10%  '

Set-TextValue 'E39' '  -3.53%  '

Set-TextValue 'D40' '17.69'
Set-TextValue 'E40' '  -2.29%  '

Set-TextValue 'D41' '0.301'
Set-TextValue 'E41' '  -4.21%  '

Set-TextValue 'D42' '4.34'
Set-TextValue 'E42' '  -5.19%  '

Set-TextValue 'D43' '1.45'
Set-TextValue 'E43' '  -7.92%  '

Set-TextValue 'D44' '1.08'
Set-TextValue 'E44' '  +1.08%  '

$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D45' '2.02'
Set-TextValue 'E45' '  -5.60%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D46' '132.42'
Set-TextValue 'E46' '  -4.15%  '

Set-TextValue 'D47' '3.31'
Set-TextValue 'E47' '  -3.68%  '

Set-TextValue 'E48' '  -2.16%  '

Set-TextValue 'D49' '0.476'
Set-TextValue 'E49' '  -6.71%  '

Set-TextValue 'D50' '0.558'
Set-TextValue 'E50' '  -2.37%  '

Set-TextValue 'D51' '0.0903'
Set-TextValue 'E51' '  -1.83%  '
